# Auto-generated Excel COM-interop script
# Applies the odds-table update for Jogos_da_Semana_FlashScore_2025-03-23.xlsx
# (see commit 'Atualizando o arquivo XLSX')

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 2.88
$ws.Range("R4").Value = 1.4
# Row 5
$ws.Range("G5").Value = 2.8
$ws.Range("I5").Value = 2.75
$ws.Range("J5").Value = 4
$ws.Range("N5").Value = 5
$ws.Range("Q5").Value = 3.5
$ws.Range("R5").Value = 1.3
$ws.Range("S5").Value = 8
$ws.Range("T5").Value = 1.08
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 2
$ws.Range("W5").Value = 2.63
$ws.Range("X5").Value = 1.44
$ws.Range("Z5").Value = 12
$ws.Range("AB5").Value = 34
$ws.Range("AF5").Value = 6
$ws.Range("AJ5").Value = 11
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 34
# Row 12
$ws.Range("G12").Value = 2.2
$ws.Range("I12").Value = 3.9
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 1.17
$ws.Range("N12").Value = 5
$ws.Range("AI12").Value = 7.5
# Row 13
$ws.Range("G13").Value = 2.55
$ws.Range("J13").Value = 3.5
$ws.Range("Z13").Value = 11
$ws.Range("AB13").Value = 26
$ws.Range("AI13").Value = 6.5
# Row 29
$ws.Range("I29").Value = 2.35
$ws.Range("K29").Value = 2.2
$ws.Range("L29").Value = 3
$ws.Range("N29").Value = 12
$ws.Range("AL29").Value = 23
# Row 30
$ws.Range("G30").Value = 1.8
$ws.Range("H30").Value = 3.7
$ws.Range("I30").Value = 4.1
$ws.Range("J30").Value = 2.4
$ws.Range("L30").Value = 4.33
$ws.Range("M30").Value = 1.04
$ws.Range("N30").Value = 13
$ws.Range("Q30").Value = 1.7
$ws.Range("R30").Value = 2.1
$ws.Range("W30").Value = 1.62
$ws.Range("X30").Value = 2.2
$ws.Range("Y30").Value = 9
$ws.Range("AF30").Value = 7
$ws.Range("AK30").Value = 13
$ws.Range("AM30").Value = 29
# Row 32
$ws.Range("G32").Value = 1.8
$ws.Range("H32").Value = 3.3
$ws.Range("L32").Value = 4.75
$ws.Range("M32").Value = 1.06
$ws.Range("N32").Value = 10
$ws.Range("O32").Value = 1.3
$ws.Range("P32").Value = 3.4
$ws.Range("Q32").Value = 2.03
$ws.Range("R32").Value = 1.83
$ws.Range("S32").Value = 3.4
$ws.Range("T32").Value = 1.3
$ws.Range("W32").Value = 1.83
$ws.Range("X32").Value = 1.83
$ws.Range("Y32").Value = 7
$ws.Range("Z32").Value = 8.5
$ws.Range("AA32").Value = 8.5
$ws.Range("AC32").Value = 15
$ws.Range("AE32").Value = 9
$ws.Range("AG32").Value = 15
$ws.Range("AI32").Value = 12
$ws.Range("AO32").Value = 301
# Row 33
$ws.Range("G33").Value = 1.33
$ws.Range("H33").Value = 4.75
$ws.Range("I33").Value = 9.5
$ws.Range("J33").Value = 1.83
$ws.Range("K33").Value = 2.38
$ws.Range("L33").Value = 8.5
$ws.Range("M33").Value = 1.06
$ws.Range("N33").Value = 10
$ws.Range("O33").Value = 1.25
$ws.Range("P33").Value = 3.75
$ws.Range("Q33").Value = 1.88
$ws.Range("R33").Value = 1.98
$ws.Range("S33").Value = 3.25
$ws.Range("T33").Value = 1.33
$ws.Range("U33").Value = 1.36
$ws.Range("V33").Value = 3
$ws.Range("W33").Value = 2.2
$ws.Range("X33").Value = 1.62
$ws.Range("Y33").Value = 6
$ws.Range("AA33").Value = 9
$ws.Range("AB33").Value = 8
$ws.Range("AC33").Value = 13
$ws.Range("AE33").Value = 10
$ws.Range("AF33").Value = 9
$ws.Range("AG33").Value = 23
$ws.Range("AJ33").Value = 41
$ws.Range("AM33").Value = 67
# Row 34
$ws.Range("G34").Value = 2.25
$ws.Range("H34").Value = 3
$ws.Range("I34").Value = 3.5
$ws.Range("J34").Value = 3.1
$ws.Range("L34").Value = 4.33
$ws.Range("M34").Value = 1.13
$ws.Range("N34").Value = 6
$ws.Range("Y34").Value = 5.5
$ws.Range("Z34").Value = 9.5
$ws.Range("AA34").Value = 10
$ws.Range("AB34").Value = 21
$ws.Range("AC34").Value = 23
$ws.Range("AG34").Value = 21
$ws.Range("AI34").Value = 7.5
$ws.Range("AJ34").Value = 15
$ws.Range("AK34").Value = 13
$ws.Range("AL34").Value = 41
$ws.Range("AM34").Value = 34
$ws.Range("AN34").Value = 51
# Row 35
$ws.Range("H35").Value = 2.9
$ws.Range("I35").Value = 2.88
$ws.Range("K35").Value = 2.1
$ws.Range("L35").Value = 3.5
$ws.Range("M35").Value = 1.06
$ws.Range("N35").Value = 10
$ws.Range("O35").Value = 1.29
$ws.Range("P35").Value = 3.5
$ws.Range("Q35").Value = 2.03
$ws.Range("R35").Value = 1.83
$ws.Range("S35").Value = 3.4
$ws.Range("T35").Value = 1.3
$ws.Range("U35").Value = 1.44
$ws.Range("V35").Value = 2.63
$ws.Range("W35").Value = 1.73
$ws.Range("X35").Value = 2
$ws.Range("Y35").Value = 8.5
$ws.Range("Z35").Value = 13
$ws.Range("AD35").Value = 29
$ws.Range("AE35").Value = 9
$ws.Range("AF35").Value = 6
$ws.Range("AH35").Value = 41
$ws.Range("AI35").Value = 9.5
$ws.Range("AK35").Value = 11
$ws.Range("AM35").Value = 23
$ws.Range("AO35").Value = 201
# Row 36
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 2.88
$ws.Range("I36").Value = 2.5
$ws.Range("L36").Value = 3.4
$ws.Range("U36").Value = 1.62
$ws.Range("V36").Value = 2.2
$ws.Range("Y36").Value = 7
$ws.Range("Z36").Value = 13
$ws.Range("AA36").Value = 12
$ws.Range("AC36").Value = 29
$ws.Range("AE36").Value = 6
$ws.Range("AH36").Value = 67
$ws.Range("AI36").Value = 6.5
$ws.Range("AJ36").Value = 11
$ws.Range("AK36").Value = 11
$ws.Range("AL36").Value = 26
$ws.Range("AM36").Value = 26
# Row 37
$ws.Range("G37").Value = 2.1
$ws.Range("I37").Value = 3.4
$ws.Range("J37").Value = 3
$ws.Range("K37").Value = 1.95
$ws.Range("N37").Value = 7.5
$ws.Range("W37").Value = 2.1
$ws.Range("X37").Value = 1.67
$ws.Range("AC37").Value = 21
# Row 38
$ws.Range("G38").Value = 1.95
$ws.Range("I38").Value = 3.8
$ws.Range("J38").Value = 2.6
$ws.Range("L38").Value = 4.33
$ws.Range("U38").Value = 1.4
$ws.Range("V38").Value = 2.75
$ws.Range("W38").Value = 1.75
$ws.Range("X38").Value = 2
$ws.Range("Y38").Value = 8
$ws.Range("Z38").Value = 9.5
$ws.Range("AD38").Value = 26
$ws.Range("AE38").Value = 11
$ws.Range("AG38").Value = 15
$ws.Range("AH38").Value = 51
$ws.Range("AI38").Value = 11
$ws.Range("AK38").Value = 13
$ws.Range("AM38").Value = 29
$ws.Range("AO38").Value = 201

Write-Output "Applied 189 odds updates across rows 4,5,12,13,29,30,32-38"
